# Adds two new "Matchday" results sheets ("8" and "9") to the workbook,
# following the same layout/style as the existing matchday sheets (1-7):
# columns B:F hold a bold/centered/bordered header row (Scorer/Team/Goals/
# Picks/Matchday), column A (bold/centered/bordered) holds a pick-count,
# and columns B:F hold the per-player results for that matchday.

$wb = $excel.ActiveWorkbook

function Add-MatchdaySheet {
    # Positional params: this COM-script engine does not bind -Name style args.
    param($wb, $SheetName, $Rows)

    # Place the new sheet right after the current last sheet, so tab order
    # stays "1", "2", ... "7", "8", "9" (matches Worksheets.Add(Before, After)).
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

    # Grab the header-row (B1:F1) and a data-cell (A2) formatting from sheet "7"
    # so the new sheet reuses the same bold/centered/bordered style (s="1").
    $formatSrc = $wb.Worksheets.Item("7")
    $formatSrc.Range("B1:F1").Copy()

    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $SheetName
    $ws.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats

    $lastRow = $Rows.Count + 1
    $formatSrc.Range("A2").Copy()
    $ws.Range("A2:A$lastRow").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Range("B1").Value = "Scorer"
    $ws.Range("C1").Value = "Team"
    $ws.Range("D1").Value = "Goals"
    $ws.Range("E1").Value = "Picks"
    $ws.Range("F1").Value = "Matchday"

    # Write column-major (all Scorer names, then all Team names, then the
    # numeric columns) so new shared-string entries land in the same order
    # the source workbook used.
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $Rows[$i]
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $r[1]
    }
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $Rows[$i]
        $row = $i + 2
        $ws.Cells.Item($row, 3).Value = $r[2]
    }
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $Rows[$i]
        $row = $i + 2
        $ws.Cells.Item($row, 4).Value = $r[3]
        $ws.Cells.Item($row, 5).Value = $r[4]
        $ws.Cells.Item($row, 6).Value = $r[5]
        $ws.Cells.Item($row, 1).Value = $r[0]
    }

    return $ws
}

$rows8 = @(
    @(55, "алькасер", "вильярреал", 0, 1, 8),
    @(76, "беседин", "динамо киев", 0, 1, 8),
    @(37, "бэйл", "тоттенхэм", 0, 8, 8),
    @(38, "винисиус", "тоттенхем", 0, 1, 8),
    @(59, "гризман", "барса", 0, 1, 8),
    @(26, "кейн", "тоттенхэм", 2, 8, 8),
    @(75, "кент", "рейнджерс", 0, 1, 8),
    @(24, "кьеза", "ювентус", 2, 1, 8),
    @(41, "марсьяль", "манчестер юнайтед", 0, 1, 8),
    @(42, "мбаппе", "псж", 1, 10, 8),
    @(72, "месси", "барселона", 1, 3, 8),
    @(65, "молина", "гранада", 1, 2, 8),
    @(36, "мората", "ювентус", 0, 1, 8),
    @(7, "морелос", "рейнджерс", 0, 2, 8),
    @(53, "морено", "вильярреал", 0, 1, 8),
    @(67, "ндомбеле", "тоттенхэм", 0, 1, 8),
    @(13, "роналду", "ювентус", 0, 12, 8),
    @(43, "рэшфорд", "манчестер юнайтед", 0, 2, 8),
    @(74, "сольдадо", "гранада", 1, 1, 8),
    @(27, "сон", "тоттенхэм", 0, 1, 8),
    @(50, "станчу", "славия", 1, 1, 8),
    @(66, "фернандеш", "манчестер юнайтед", 0, 3, 8),
    @(20, "холанд", "боруссия д", 2, 13, 8),
)

$rows9 = @(
    @(0, "аспас", "сельта", 0, 4, 9),
    @(62, "бейл", "тоттенхэм", 0, 2, 9),
    @(36, "бен-йеддер", "монако", 0, 3, 9),
    @(46, "беседин", "динамо киев", 1, 1, 9),
    @(48, "броя", "витесс", 0, 1, 9),
    @(2, "буяльский", "динамо киев", 0, 1, 9),
    @(73, "вамангитука", "штутгарт", 0, 1, 9),
    @(30, "ван де стрек", "утрехт", 0, 4, 9),
    @(59, "гладкий", "заря", 0, 1, 9),
    @(11, "де паул", "удинезе", 1, 1, 9),
    @(34, "дестро", "дженоа", 0, 1, 9),
    @(32, "захеди", "заря", 0, 1, 9),
    @(3, "ибрагимович", "милан", 0, 1, 9),
    @(45, "ингс", "саутгемптон", 0, 1, 9),
    @(27, "инсинье", "наполи", 0, 1, 9),
    @(4, "каладжич", "штутгарт", 1, 5, 9),
    @(29, "кейн", "тоттенхэм", 0, 9, 9),
    @(37, "кессье", "милан", 0, 3, 9),
    @(72, "крамарич", "хоффенхайм", 0, 1, 9),
    @(67, "кулибали", "штутгарт", 0, 1, 9),
    @(55, "ларссон", "спартак", 0, 3, 9),
    @(35, "леау", "милан", 0, 1, 9),
    @(13, "льоренте", "удинезе", 0, 1, 9),
    @(54, "ляказетт", "арсенал", 1, 2, 9),
    @(53, "махи", "утрехт", 0, 2, 9),
    @(33, "мина", "сельта", 0, 2, 9),
    @(7, "мопе", "брайтон", 0, 1, 9),
    @(66, "обамеянг", "арсенал", 0, 2, 9),
    @(26, "опенда", "витесс", 0, 1, 9),
    @(16, "понсе", "спартак", 0, 1, 9),
    @(8, "промес", "спартак", 0, 3, 9),
    @(10, "соболев", "спартак", 1, 2, 9),
    @(63, "уорд-проуз", "саутгемптон", 0, 1, 9),
    @(61, "фолланд", "монако", 0, 2, 9),
    @(1, "фомин", "динамо москва", 0, 2, 9),
    @(65, "цыганков", "динамо киев", 0, 6, 9),
    @(76, "эдегор", "арсенал", 1, 1, 9),
    @(50, "эндо", "штутгарт", 0, 2, 9),
)

Add-MatchdaySheet $wb "8" $rows8 | Out-Null
Add-MatchdaySheet $wb "9" $rows9 | Out-Null
